$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# Header date line
Replace-Text "2026-01-22 Thursday" "2026-01-23 Friday"

# Multiplication problems, in document order (top-to-bottom, left-to-right)
Replace-Text "781×6=" "113×6="
Replace-Text "214×8=" "742×6="
Replace-Text "842×5=" "551×3="
Replace-Text "473×3=" "559×4="
Replace-Text "787×8=" "373×9="
Replace-Text "916×9=" "297×2="
Replace-Text "479×8=" "837×9="
Replace-Text "179×5=" "658×8="
Replace-Text "141×7=" "399×5="
Replace-Text "578×3=" "705×8="
Replace-Text "307×8=" "555×5="
Replace-Text "512×6=" "152×9="
Replace-Text "499×4=" "621×6="
Replace-Text "870×4=" "249×8="
Replace-Text "903×2=" "826×7="
Replace-Text "756×2=" "121×8="
Replace-Text "749×7=" "314×2="
Replace-Text "675×4=" "784×3="
Replace-Text "610×4=" "499×4="
Replace-Text "977×7=" "858×5="
Replace-Text "438×2=" "638×2="
Replace-Text "785×2=" "632×4="
Replace-Text "418×5=" "285×6="
Replace-Text "714×6=" "238×7="
Replace-Text "361×7=" "440×4="
